$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain numeric" string (e.g. "594.17") need to be
# protected from Excel auto-converting the text into a Number: briefly mark the
# cell as Text, write the value, then restore the Normal style so no stray
# number-format style lingers on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '63.256.41'
$ws.Range('E2').Value = '  -1.26%  '

$ws.Range('D3').Value = '3.227.69'
$ws.Range('E3').Value = '  +2.38%  '

$ws.Range('E4').Value = '  -0.10%  '

Set-TextValue $ws.Range('D5') '594.17'
$ws.Range('E5').Value = '  -1.30%  '

Set-TextValue $ws.Range('D6') '141.30'
$ws.Range('E6').Value = '  -1.50%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').Value = '3.223.62'
$ws.Range('E8').Value = '  +2.63%  '

Set-TextValue $ws.Range('D9') '0.519'
$ws.Range('E9').Value = '  -1.80%  '

Set-TextValue $ws.Range('D10') '0.147'
$ws.Range('E10').Value = '  -1.66%  '

Set-TextValue $ws.Range('D11') '5.35'
$ws.Range('E11').Value = '  -0.80%  '

Set-TextValue $ws.Range('D12') '0.464'
$ws.Range('E12').Value = '  -0.92%  '

Set-TextValue $ws.Range('D13') '0.0000246'
$ws.Range('E13').Value = '  -3.21%  '

Set-TextValue $ws.Range('D14') '34.33'
$ws.Range('E14').Value = '  -2.26%  '

$ws.Range('D15').Value = '3.760.40'
$ws.Range('E15').Value = '  +2.54%  '

$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('D17').Value = '3.238.24'
$ws.Range('E17').Value = '  +3.09%  '

$ws.Range('D18').Value = '63.272.18'
$ws.Range('E18').Value = '  -1.29%  '

Set-TextValue $ws.Range('D19') '6.77'
$ws.Range('E19').Value = '  -1.72%  '

Set-TextValue $ws.Range('D20') '474.13'
$ws.Range('E20').Value = '  -3.15%  '

Set-TextValue $ws.Range('D21') '14.14'
$ws.Range('E21').Value = '  -4.02%  '

Set-TextValue $ws.Range('D22') '0.725'
$ws.Range('E22').Value = '  +1.50%  '

Set-TextValue $ws.Range('D23') '7.91'
$ws.Range('E23').Value = '  +1.85%  '

Set-TextValue $ws.Range('D24') '84.03'
$ws.Range('E24').Value = '  -4.81%  '

Set-TextValue $ws.Range('D25') '13.13'
$ws.Range('E25').Value = '  -1.52%  '

$ws.Range('E26').Value = '  -0.13%  '

Set-TextValue $ws.Range('D27') '7.51'
$ws.Range('E27').Value = '  +6.78%  '

Set-TextValue $ws.Range('D28') '2.73'
$ws.Range('E28').Value = '  -1.29%  '

Set-TextValue $ws.Range('D29') '8.08'
$ws.Range('E29').Value = '  -1.93%  '

$ws.Range('E30').Value = '  +1.24%  '

Set-TextValue $ws.Range('D31') '27.36'
$ws.Range('E31').Value = '  -1.30%  '

$ws.Range('E32').Value = '  -0.03%  '

Set-TextValue $ws.Range('D33') '0.107'
$ws.Range('E33').Value = '  -4.30%  '

Set-TextValue $ws.Range('D34') '2.52'
$ws.Range('E34').Value = '  -5.18%  '

$ws.Range('E35').Value = '  -2.13%  '

Set-TextValue $ws.Range('D36') '5.90'
$ws.Range('E36').Value = '  -2.92%  '

$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('D38').Value = '0.0₃0706'
$ws.Range('E38').Value = '  -5.79%  '

Set-TextValue $ws.Range('D39') '0.0392'
$ws.Range('E39').Value = '  -1.53%  '

Set-TextValue $ws.Range('D40') '422.23'
$ws.Range('E40').Value = '  -2.84%  '

Set-TextValue $ws.Range('D41') '8.37'
$ws.Range('E41').Value = '  -0.21%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.973.55'
$ws.Range('E42').Value = '  +1.30%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D43') '2.74'
$ws.Range('E43').Value = '  -7.92%  '

Set-TextValue $ws.Range('D44') '0.109'
$ws.Range('E44').Value = '  -9.04%  '

Set-TextValue $ws.Range('D45') '0.266'
$ws.Range('E45').Value = '  +1.73%  '

$ws.Range('E46').Value = '  -1.96%  '

Set-TextValue $ws.Range('D47') '2.36'
$ws.Range('E47').Value = '  -2.18%  '

$ws.Range('E48').Value = '  +0.03%  '

Set-TextValue $ws.Range('D49') '25.88'
$ws.Range('E49').Value = '  -0.31%  '

$ws.Range('E50').Value = '  -0.78%  '

Set-TextValue $ws.Range('D51') '120.96'
$ws.Range('E51').Value = '  +0.45%  '
